# 10.1.1 SDG indicator sheet -- reshuffle the trilingual title/subtitle rows,
# add the "(в процентах)/(in percent)" subtitle text, turn the separator row
# into a styled thin row, re-center the year-header row vertically, and add
# a new 2023 data column (L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 - big trilingual title. Columns get reordered: Kyrgyz first, then
# Russian, then English (previously Russian, English, Kyrgyz).
# ---------------------------------------------------------------------
$ws.Range("A1").Value2 = " 10.1.1 40 пайыздан аз камсыздалган калктын жана жалпы калктын арасындагы  үй чарбалардын кирешелеринин  калктын жан башына алгандагы өсүү темпи "
$ws.Range("B1").Value2 = "10.1.1 Темпы роста доходов домохозяйств  на душу населения среди наименее обеспеченных 40 процентов населения и среди населения в целом"
$ws.Range("C1").Value2 = "10.1.1 Growth rates of household expenditure or income per capita among the bottom 40 per cent of the population and the total population"

# C1 used to be indented; the new layout drops that indent.
$ws.Range("C1").IndentLevel = 0

# ---------------------------------------------------------------------
# Row 2 - "(percent)" subtitle, same reorder plus new wording in RU/EN.
# ---------------------------------------------------------------------
$ws.Range("A2").Value2 = "(пайыз менен)"
$ws.Range("B2").Value2 = "(в процентах)"
$ws.Range("C2").Value2 = "(in percent)"

# Give the rest of row 1 / row 2 (D:K) the same vertically-centered
# formatting so the row reads as one visual band.
$ws.Range("D1:K2").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Row 3 - thin separator row: taller (13.5 instead of 10.5) and now
# carries real (blank) cells so the bottom border shows across A:K.
# ---------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 13.5
$ws.Range("A3:K3").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Row 4 - header row (labels + years): switch vertical alignment from
# top to center, and add the 2023 column header.
# ---------------------------------------------------------------------
$ws.Range("A4:L4").VerticalAlignment = -4108
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value2 = 2023

# ---------------------------------------------------------------------
# Row 5 / 6 - data rows: add the 2023 values, copying the format from
# the last existing column so the number format / borders match.
# ---------------------------------------------------------------------
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value2 = 4.9000000000000004

$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value2 = 4.5999999999999996

$ws.Rows.Item(6).RowHeight = 28.5

# ---------------------------------------------------------------------
# Columns A:C all become a uniform 41-wide (previously 43 / 43.86 / 41.86).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 40.166666666666664
$ws.Columns.Item(2).ColumnWidth = 40.166666666666664
$ws.Columns.Item(3).ColumnWidth = 40.166666666666664

$ws.Range("A1").Select()
